$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.168.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "'1.825.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'235.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").Value = "'0.6010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.03%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.07074"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.78%  "

$ws.Range("D9").Value = "'0.2791"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.58%  "

$ws.Range("D10").Value = "'23.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.20%  "

$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").Value = "'1.827.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "'4.791"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.23%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6301"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.51%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.000009951"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("D16").Value = "'2.073.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "'78.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.56%  "

$ws.Range("D18").Value = "'5.857"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.07%  "

$ws.Range("D19").Value = "'29.158.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").Value = "'226.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("E22").Value = "  -4.68%  "

$ws.Range("D23").Value = "'6.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.75%  "

$ws.Range("D24").Value = "'0.9994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'154.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.96%  "

$ws.Range("D26").Value = "'8.025"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.32%  "

$ws.Range("D27").Value = "'0.1299"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.49%  "

$ws.Range("E28").Value = "  -4.52%  "

$ws.Range("D29").Value = "'1.492"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("D30").Value = "'0.06217"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -15.11%  "

$ws.Range("D31").Value = "'1.450"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "

$ws.Range("D32").Value = "'3.832"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.07%  "

$ws.Range("D33").Value = "'3.799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.83%  "

$ws.Range("D34").Value = "'1.122"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.43%  "

$ws.Range("D35").Value = "'1.742"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("D36").Value = "'0.6405"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.04%  "

$ws.Range("E37").Value = "  -1.60%  "

$ws.Range("D38").Value = "'1.213.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("D39").Value = "'2.730"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.05%  "

$ws.Range("D40").Value = "'0.01733"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.43%  "

$ws.Range("D41").Value = "'6.492"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.51%  "

$ws.Range("D42").Value = "'0.9054"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.83%  "

$ws.Range("D44").Value = "'1.982.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").Value = "'100.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").Value = "'62.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.02%  "

$ws.Range("D47").Value = "'0.00000000118"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").Value = "'8.520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("E49").Value = "  -6.39%  "

$ws.Range("D50").Value = "'0.4556"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").Value = "'0.05501"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.67%  "
